$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H3").Value = "Resolved"
$ws.Range("H4").Value = "Resolved"

$ws.Activate()
$ws.Range("H4").Select()
